$d = $word.ActiveDocument

$pairs = @(
    @{old="780×9=7020"; new="372×2=744"},
    @{old="873×6=5238"; new="735×6=4410"},
    @{old="502×9=4518"; new="121×6=726"},
    @{old="938×3=2814"; new="354×2=708"},
    @{old="613×9=5517"; new="854×7=5978"},
    @{old="834×6=5004"; new="335×8=2680"},
    @{old="333×9=2997"; new="817×8=6536"},
    @{old="566×9=5094"; new="313×2=626"},
    @{old="936×6=5616"; new="627×2=1254"},
    @{old="439×5=2195"; new="804×6=4824"},
    @{old="962×8=7696"; new="296×9=2664"},
    @{old="904×8=7232"; new="814×9=7326"},
    @{old="840×7=5880"; new="175×9=1575"},
    @{old="995×6=5970"; new="285×4=1140"},
    @{old="505×6=3030"; new="226×2=452"},
    @{old="962×9=8658"; new="944×2=1888"},
    @{old="746×3=2238"; new="933×7=6531"},
    @{old="877×7=6139"; new="128×8=1024"},
    @{old="556×4=2224"; new="922×9=8298"},
    @{old="475×8=3800"; new="697×2=1394"},
    @{old="211×3=633"; new="516×5=2580"},
    @{old="618×4=2472"; new="442×3=1326"},
    @{old="756×6=4536"; new="319×4=1276"},
    @{old="853×3=2559"; new="610×5=3050"},
    @{old="816×3=2448"; new="865×3=2595"}
)

foreach ($p in $pairs) {
    $d.Content.Find.Execute($p.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $p.new, 2)
}
